# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.179.62"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "2.589.72"
$ws.Range("E3").Value = "  +8.46%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'305.29"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").Value = "'99.42"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").Value = "'0.598"
$ws.Range("E7").Value = "  +5.35%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.568"
$ws.Range("E9").Value = "  +11.58%  "

$ws.Range("D10").Value = "'38.41"
$ws.Range("E10").Value = "  +11.29%  "

$ws.Range("D11").Value = "'0.0831"
$ws.Range("E11").Value = "  +5.12%  "

$ws.Range("D12").Value = "'8.06"
$ws.Range("E12").Value = "  +12.93%  "

$ws.Range("D13").Value = "2.992.00"
$ws.Range("E13").Value = "  +8.89%  "

$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("D15").Value = "2.620.00"
$ws.Range("E15").Value = "  +9.83%  "

$ws.Range("D16").Value = "'0.891"
$ws.Range("E16").Value = "  +8.90%  "

$ws.Range("D17").Value = "'14.78"
$ws.Range("E17").Value = "  +7.72%  "

$ws.Range("D18").Value = "46.287.05"
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  +2.76%  "

$ws.Range("D20").Value = "0.0₃0996"
$ws.Range("E20").Value = "  +4.42%  "

$ws.Range("D21").Value = "'6.62"
$ws.Range("E21").Value = "  +9.48%  "

$ws.Range("D22").Value = "'70.81"
$ws.Range("E22").Value = "  +5.62%  "

$ws.Range("D23").Value = "'253.62"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").Value = "'2.98"
$ws.Range("E24").Value = "  +6.64%  "

$ws.Range("E25").Value = "  +15.19%  "

$ws.Range("D26").Value = "'27.66"
$ws.Range("E26").Value = "  +31.04%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").Value = "'10.37"
$ws.Range("E28").Value = "  +6.21%  "

$ws.Range("D29").Value = "'39.64"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  +3.24%  "

$ws.Range("D31").Value = "'6.11"
$ws.Range("E31").Value = "  +10.31%  "

$ws.Range("D32").Value = "'3.67"
$ws.Range("E32").Value = "  -3.95%  "

$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.93"
$ws.Range("E33").Value = "  +4.32%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'2.29"
$ws.Range("E34").Value = "  +17.88%  "

$ws.Range("D35").Value = "'151.90"
$ws.Range("E35").Value = "  +3.43%  "

$ws.Range("D36").Value = "'0.0827"
$ws.Range("E36").Value = "  +7.08%  "

$ws.Range("D37").Value = "'0.116"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").Value = "'0.121"
$ws.Range("E38").Value = "  +4.77%  "

$ws.Range("D39").Value = "'4.16"
$ws.Range("E39").Value = "  +6.39%  "

$ws.Range("D40").Value = "'15.48"
$ws.Range("E40").Value = "  +4.28%  "

$ws.Range("D41").Value = "'3.57"
$ws.Range("E41").Value = "  +10.62%  "

$ws.Range("D42").Value = "'0.0320"
$ws.Range("E42").Value = "  +6.85%  "

$ws.Range("D43").Value = "2.041.85"
$ws.Range("E43").Value = "  +5.21%  "

$ws.Range("D44").Value = "'18.93"
$ws.Range("E44").Value = "  +33.27%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("D46").Value = "'90.62"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.11"
$ws.Range("E47").Value = "  +7.48%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'108.76"
$ws.Range("E48").Value = "  +10.25%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.77"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("D50").Value = "2.843.65"
$ws.Range("E50").Value = "  +8.54%  "

$ws.Range("D51").Value = "'0.198"
$ws.Range("E51").Value = "  +6.82%  "
